$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "243.02"
$cell.Style = "Normal"

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "23.71"
$cell.Style = "Normal"

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "5.227"
$cell.Style = "Normal"

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "0.05757"
$cell.Style = "Normal"

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "6.405"
$cell.Style = "Normal"

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "3.228"
$cell.Style = "Normal"

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.8057"
$cell.Style = "Normal"

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.8789"
$cell.Style = "Normal"

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.07076"
$cell.Style = "Normal"

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.03148"
$cell.Style = "Normal"

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.03040"
$cell.Style = "Normal"

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "0.09314"
$cell.Style = "Normal"

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "3.814"
$cell.Style = "Normal"

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "0.001543"
$cell.Style = "Normal"

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.04712"
$cell.Style = "Normal"

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.0006022"
$cell.Style = "Normal"

$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = "17OneONEWorstin24h"
$cell.Style = "Normal"

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.006264"
$cell.Style = "Normal"

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "0.004054"
$cell.Style = "Normal"

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.00008706"
$cell.Style = "Normal"

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "3.544"
$cell.Style = "Normal"

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "2.150"
$cell.Style = "Normal"

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "0.3162"
$cell.Style = "Normal"

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "0.0002329"
$cell.Style = "Normal"

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.03722"
$cell.Style = "Normal"

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.006275"
$cell.Style = "Normal"

$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = "40KickTokenKICK"
$cell.Style = "Normal"

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.1044"
$cell.Style = "Normal"

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.002472"
$cell.Style = "Normal"

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.007146"
$cell.Style = "Normal"

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.00005335"
$cell.Style = "Normal"

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.00000000750"
$cell.Style = "Normal"

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.5352"
$cell.Style = "Normal"

$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = "46CoinbaseStockTokenCOIN"
$cell.Style = "Normal"

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "0.002361"
$cell.Style = "Normal"

$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = "47BOLOBOLOBestin24h"
$cell.Style = "Normal"

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.00002101"
$cell.Style = "Normal"

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.0002001"
$cell.Style = "Normal"
